$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly.
$ws.Range('D2').Value = '30.129.44'
$ws.Range('D3').Value = '1.858.03'
$ws.Range('E3').Value = '  -2.83%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -1.79%  '
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('E10').Value = '  +3.02%  '
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('E12').Value = '  -5.90%  '
$ws.Range('D13').Value = '1.855.47'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('E16').Value = '  +3.05%  '
$ws.Range('D17').Value = '30.154.87'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '2.102.50'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('E22').Value = '  -3.50%  '
$ws.Range('E24').Value = '  -2.45%  '
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('E28').Value = '  -7.37%  '
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('E30').Value = '  -3.46%  '
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('E33').Value = '  -3.59%  '
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  -4.12%  '
$ws.Range('E36').Value = '  -2.38%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('E42').Value = '  -3.23%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '1.031.53'
$ws.Range('E48').Value = '  +8.13%  '
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('E50').Value = '  +4.26%  '
$ws.Range('E51').Value = '  -2.36%  '

# Numeric-looking text values (e.g. "1.0000", "0.9999"): Excel would
# auto-convert these to numbers on plain assignment, losing the exact
# text representation (trailing zeros, etc). Force text entry by setting
# the cell to Text format, assign, then restore the Normal style so the
# cell keeps its original (default) style index.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4695'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2819'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06550'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07800'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.088'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6717'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.86'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9999'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.451'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007235'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '167.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.313'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.930'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.342'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09646'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.408'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.470'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.105'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04683'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6981'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.091'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9992'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.703'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01864'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.325'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8618'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.940'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4168'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.254'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.108'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.82'
$ws.Range('D51').Style = 'Normal'
